$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_3_9_10"
$ws.Range("B2").Value = 0.1984207691934646
$ws.Range("C2").Value = -0.755002438371885
$ws.Range("D2").Value = 0.2729627960098634
$ws.Range("E2").Value = -0.3051436890544714
$ws.Range("F2").Value = 0.8871120810508728
$ws.Range("G2").Value = 3.326462030410767
$ws.Range("H2").Value = 1.193012475967407
$ws.Range("I2").Value = 2.322487592697144

$ws.Range("A3").Value = "model_3_9_11"
$ws.Range("B3").Value = 0.2895599609140579
$ws.Range("C3").Value = -0.4545896609217528
$ws.Range("D3").Value = 0.3121754711659134
$ws.Range("E3").Value = -0.1187239159318012
$ws.Range("F3").Value = 0.7862477898597717
$ws.Range("G3").Value = 2.757054328918457
$ws.Range("H3").Value = 1.12866735458374
$ws.Range("I3").Value = 1.990755915641785

$ws.Range("A4").Value = "model_3_9_0"
$ws.Range("B4").Value = 0.3145054881783034
$ws.Range("C4").Value = 0.3085973256153984
$ws.Range("D4").Value = 0.3038203107847639
$ws.Range("E4").Value = 0.3080131231787764
$ws.Range("F4").Value = 0.7586404085159302
$ws.Range("G4").Value = 1.31049656867981
$ws.Range("H4").Value = 1.142377614974976
$ws.Range("I4").Value = 1.231382250785828

$ws.Range("A5").Value = "model_3_9_1"
$ws.Range("B5").Value = 0.3550102809061463
$ws.Range("C5").Value = 0.3579023695474768
$ws.Range("D5").Value = 0.3501438402067579
$ws.Range("E5").Value = 0.3559184022467831
$ws.Range("F5").Value = 0.7138136029243469
$ws.Range("G5").Value = 1.217043042182922
$ws.Range("H5").Value = 1.066364288330078
$ws.Range("I5").Value = 1.146135449409485

$ws.Range("A6").Value = "model_3_9_24"
$ws.Range("B6").Value = 0.3621708436784778
$ws.Range("C6").Value = -0.6376081593365817
$ws.Range("D6").Value = 0.5032803830974275
$ws.Range("E6").Value = -0.138998731964453
$ws.Range("F6").Value = 0.7058889269828796
$ws.Range("G6").Value = 3.103950977325439
$ws.Range("H6").Value = 0.8150789141654968
$ws.Range("I6").Value = 2.026834487915039

$ws.Range("A7").Value = "model_3_9_2"
$ws.Range("B7").Value = 0.369967805265397
$ws.Range("C7").Value = 0.3431263832225573
$ws.Range("D7").Value = 0.3784496660375636
$ws.Range("E7").Value = 0.3598701885935421
$ws.Range("F7").Value = 0.6972599625587463
$ws.Range("G7").Value = 1.245049715042114
$ws.Range("H7").Value = 1.019916534423828
$ws.Range("I7").Value = 1.139103412628174

$ws.Range("A8").Value = "model_3_9_6"
$ws.Range("B8").Value = 0.387306005258068
$ws.Range("C8").Value = 0.3759687083996167
$ws.Range("D8").Value = 0.3029261734684242
$ws.Range("E8").Value = 0.345617596094594
$ws.Range("F8").Value = 0.6780716776847839
$ws.Range("G8").Value = 1.182799816131592
$ws.Range("H8").Value = 1.143844842910767
$ws.Range("I8").Value = 1.164465665817261

$ws.Range("A9").Value = "model_3_9_12"
$ws.Range("B9").Value = 0.3880282884643048
$ws.Range("C9").Value = -0.1059053172055744
$ws.Range("D9").Value = 0.3089604553212345
$ws.Range("E9").Value = 0.07650559093311216
$ws.Range("F9").Value = 0.6772723197937012
$ws.Range("G9").Value = 2.096152067184448
$ws.Range("H9").Value = 1.1339430809021
$ws.Range("I9").Value = 1.643347144126892

$ws.Range("A10").Value = "model_3_9_7"
$ws.Range("B10").Value = 0.3935539720499002
$ws.Range("C10").Value = 0.4227697449421465
$ws.Range("D10").Value = 0.2405798742345481
$ws.Range("E10").Value = 0.3449519582033592
$ws.Range("F10").Value = 0.6711570024490356
$ws.Range("G10").Value = 1.0940922498703
$ws.Range("H10").Value = 1.246150374412537
$ws.Range("I10").Value = 1.165650129318237

$ws.Range("A11").Value = "model_3_9_13"
$ws.Range("B11").Value = 0.3938834277856962
$ws.Range("C11").Value = -0.2161007450367096
$ws.Range("D11").Value = 0.3741286487385779
$ws.Range("E11").Value = 0.04264413064367489
$ws.Range("F11").Value = 0.6707924604415894
$ws.Range("G11").Value = 2.305018424987793
$ws.Range("H11").Value = 1.027006983757019
$ws.Range("I11").Value = 1.703603029251099

$ws.Range("A12").Value = "model_3_9_9"
$ws.Range("B12").Value = 0.3981693583743471
$ws.Range("C12").Value = -0.01563597325747801
$ws.Range("D12").Value = 0.349173321943922
$ws.Range("E12").Value = 0.1448582725017978
$ws.Range("F12").Value = 0.6660491228103638
$ws.Range("G12").Value = 1.925053954124451
$ws.Range("H12").Value = 1.067956805229187
$ws.Range("I12").Value = 1.521714329719543

$ws.Range("A13").Value = "model_3_9_23"
$ws.Range("B13").Value = 0.402669787545696
$ws.Range("C13").Value = -0.4712773376827364
$ws.Range("D13").Value = 0.4902841337549191
$ws.Range("E13").Value = -0.0508444720572141
$ws.Range("F13").Value = 0.6610685586929321
$ws.Range("G13").Value = 2.788684368133545
$ws.Range("H13").Value = 0.8364047408103943
$ws.Range("I13").Value = 1.869965076446533

$ws.Range("A14").Value = "model_3_9_3"
$ws.Range("B14").Value = 0.4035145122776048
$ws.Range("C14").Value = 0.3992631471065085
$ws.Range("D14").Value = 0.4172736054200461
$ws.Range("E14").Value = 0.4083728578798319
$ws.Range("F14").Value = 0.6601336598396301
$ws.Range("G14").Value = 1.138646960258484
$ws.Range("H14").Value = 0.9562094211578369
$ws.Range("I14").Value = 1.052793383598328

$ws.Range("A15").Value = "model_3_9_8"
$ws.Range("B15").Value = 0.404129692940641
$ws.Range("C15").Value = 0.01502780572045603
$ws.Range("D15").Value = 0.3558815368078704
$ws.Range("E15").Value = 0.1650609415119059
$ws.Range("F15").Value = 0.6594528555870056
$ws.Range("G15").Value = 1.866933345794678
$ws.Range("H15").Value = 1.056949138641357
$ws.Range("I15").Value = 1.485764026641846

$ws.Range("A16").Value = "model_3_9_4"
$ws.Range("B16").Value = 0.4046620651602759
$ws.Range("C16").Value = 0.4260307037253668
$ws.Range("D16").Value = 0.3672791853375699
$ws.Range("E16").Value = 0.4017725780078807
$ws.Range("F16").Value = 0.6588636636734009
$ws.Range("G16").Value = 1.087911248207092
$ws.Range("H16").Value = 1.038246512413025
$ws.Range("I16").Value = 1.064538478851318

$ws.Range("A17").Value = "model_3_9_5"
$ws.Range("B17").Value = 0.4112000840733212
$ws.Range("C17").Value = 0.4389838797550596
$ws.Range("D17").Value = 0.3677690165196428
$ws.Range("E17").Value = 0.4092890438772405
$ws.Range("F17").Value = 0.6516279578208923
$ws.Range("G17").Value = 1.06335973739624
$ws.Range("H17").Value = 1.037442684173584
$ws.Range("I17").Value = 1.051162958145142

$ws.Range("A18").Value = "model_3_9_15"
$ws.Range("B18").Value = 0.4484990523576957
$ws.Range("C18").Value = -0.2126454687806933
$ws.Range("D18").Value = 0.5252841202625099
$ws.Range("E18").Value = 0.1101868670557677
$ws.Range("F18").Value = 0.6103490591049194
$ws.Range("G18").Value = 2.298469066619873
$ws.Range("H18").Value = 0.7789725065231323
$ws.Range("I18").Value = 1.583411812782288

$ws.Range("A19").Value = "model_3_9_21"
$ws.Range("B19").Value = 0.4520855742660452
$ws.Range("C19").Value = -0.2721017681446136
$ws.Range("D19").Value = 0.4931456847323811
$ws.Range("E19").Value = 0.06271406373543065
$ws.Range("F19").Value = 0.6063798666000366
$ws.Range("G19").Value = 2.411163806915283
$ws.Range("H19").Value = 0.8317091464996338
$ws.Range("I19").Value = 1.66788911819458

$ws.Range("A20").Value = "model_3_9_20"
$ws.Range("B20").Value = 0.453338803427689
$ws.Range("C20").Value = -0.2663311599516927
$ws.Range("D20").Value = 0.4954020795959634
$ws.Range("E20").Value = 0.06694628741091257
$ws.Range("F20").Value = 0.6049928665161133
$ws.Range("G20").Value = 2.400225877761841
$ws.Range("H20").Value = 0.8280065655708313
$ws.Range("I20").Value = 1.66035783290863

$ws.Range("A21").Value = "model_3_9_22"
$ws.Range("B21").Value = 0.4563459571299018
$ws.Range("C21").Value = -0.2627906434888923
$ws.Range("D21").Value = 0.5006146442718449
$ws.Range("E21").Value = 0.07120475685002259
$ws.Range("F21").Value = 0.6016648411750793
$ws.Range("G21").Value = 2.393515348434448
$ws.Range("H21").Value = 0.8194531202316284
$ws.Range("I21").Value = 1.652779817581177

$ws.Range("A22").Value = "model_3_9_14"
$ws.Range("B22").Value = 0.4932895077514599
$ws.Range("C22").Value = -0.01987590955977936
$ws.Range("D22").Value = 0.5066093156244913
$ws.Range("E22").Value = 0.2107850791132883
$ws.Range("F22").Value = 0.5607792139053345
$ws.Range("G22").Value = 1.933090448379517
$ws.Range("H22").Value = 0.8096163272857666
$ws.Range("I22").Value = 1.404398441314697

$ws.Range("A23").Value = "model_3_9_16"
$ws.Range("B23").Value = 0.5369693320818741
$ws.Range("C23").Value = 0.08584842525831116
$ws.Range("D23").Value = 0.5265314321533905
$ws.Range("E23").Value = 0.2790495453343343
$ws.Range("F23").Value = 0.5124384760856628
$ws.Range("G23").Value = 1.732698678970337
$ws.Range("H23").Value = 0.7769256830215454
$ws.Range("I23").Value = 1.282922744750977

$ws.Range("A24").Value = "model_3_9_19"
$ws.Range("B24").Value = 0.5388831185299665
$ws.Range("C24").Value = 0.02374760375531448
$ws.Range("D24").Value = 0.5362127203631808
$ws.Range("E24").Value = 0.2482327588245357
$ws.Range("F24").Value = 0.5103205442428589
$ws.Range("G24").Value = 1.85040557384491
$ws.Range("H24").Value = 0.7610394358634949
$ws.Range("I24").Value = 1.3377605676651

$ws.Range("A25").Value = "model_3_9_17"
$ws.Range("B25").Value = 0.542203478839279
$ws.Range("C25").Value = 0.06617840093496341
$ws.Range("D25").Value = 0.544506536600988
$ws.Range("E25").Value = 0.2757573827257228
$ws.Range("F25").Value = 0.5066457986831665
$ws.Range("G25").Value = 1.769981622695923
$ws.Range("H25").Value = 0.7474299073219299
$ws.Range("I25").Value = 1.288781046867371

$ws.Range("A26").Value = "model_3_9_18"
$ws.Range("B26").Value = 0.5616540272509607
$ws.Range("C26").Value = 0.1080373977190219
$ws.Range("D26").Value = 0.5561638952554557
$ws.Range("E26").Value = 0.3044196478657872
$ws.Range("F26").Value = 0.4851198494434357
$ws.Range("G26").Value = 1.690641403198242
$ws.Range("H26").Value = 0.7283011078834534
$ws.Range("I26").Value = 1.237776756286621
